$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111896637
$ws.Range("B2").Value = 90466
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 4769
$ws.Range("F2").Value = "Svavelriska"
$ws.Range("G2").Value = "Lactarius scrobiculatus"
$ws.Range("H2").Value = "(Scop.:Fr.) Fr."
$ws.Range("P2").Value = "Kratte masugn, Gstr"
$ws.Range("Q2").Value = 575088
$ws.Range("R2").Value = 6703396
$ws.Range("AW2").Value = "Philipp Weiss"
$ws.Range("AX2").Value = "Philipp Weiss"

# Row 3
$ws.Range("A3").Value = 111896640
$ws.Range("B3").Value = 90466
$ws.Range("Q3").Value = 575025
$ws.Range("R3").Value = 6703369

# Row 4
$ws.Range("A4").Value = 111896638
$ws.Range("B4").Value = 90466
$ws.Range("Q4").Value = 575087
$ws.Range("R4").Value = 6703393

# Row 5
$ws.Range("A5").Value = 111896654
$ws.Range("B5").Value = 89317
$ws.Range("E5").Value = 3215
$ws.Range("F5").Value = "Rödgul trumpetsvamp"
$ws.Range("G5").Value = "Craterellus lutescens"
$ws.Range("H5").Value = "(Fr.) Fr."
$ws.Range("Q5").Value = 575073
$ws.Range("R5").Value = 6703422

# Row 6
$ws.Range("A6").Value = 111884133
$ws.Range("B6").Value = 89033
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 3286
$ws.Range("F6").Value = "Flattoppad klubbsvamp"
$ws.Range("G6").Value = "Clavariadelphus truncatus"
$ws.Range("H6").Value = "(Quél.) Donk"
$ws.Range("Q6").Value = 575059
$ws.Range("R6").Value = 6703389

# Row 7
$ws.Range("A7").Value = 111896643
$ws.Range("B7").Value = 90466
$ws.Range("Q7").Value = 575039
$ws.Range("R7").Value = 6703416

# Row 8
$ws.Range("A8").Value = 111884093
$ws.Range("B8").Value = 98961
$ws.Range("E8").Value = 222498
$ws.Range("F8").Value = "Blåsippa"
$ws.Range("G8").Value = "Hepatica nobilis"
$ws.Range("H8").Value = "Schreb."
$ws.Range("P8").Value = "Kopparåsen (Kopparåsen), Gstr"
$ws.Range("Q8").Value = 575066
$ws.Range("R8").Value = 6703388
$ws.Range("AW8").Value = "Patric Engfeldt"
$ws.Range("AX8").Value = "Patric Engfeldt"

# Row 9
$ws.Range("A9").Value = 111896690
$ws.Range("B9").Value = 90821
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 5964
$ws.Range("F9").Value = "Fjällig taggsvamp s.str."
$ws.Range("G9").Value = "Sarcodon imbricatus s.str."
$ws.Range("H9").Value = "(L.:Fr.) P.Karst."
$ws.Range("P9").Value = "Kratte masugn, Gstr"
$ws.Range("Q9").Value = 575060
$ws.Range("R9").Value = 6703377
$ws.Range("AW9").Value = "Philipp Weiss"
$ws.Range("AX9").Value = "Philipp Weiss"

# Row 10
$ws.Range("A10").Value = 111896639
$ws.Range("B10").Value = 90466
$ws.Range("Q10").Value = 575089
$ws.Range("R10").Value = 6703380

# Row 11
$ws.Range("A11").Value = 111896644
$ws.Range("B11").Value = 90466
$ws.Range("E11").Value = 4769
$ws.Range("F11").Value = "Svavelriska"
$ws.Range("G11").Value = "Lactarius scrobiculatus"
$ws.Range("H11").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q11").Value = 575036
$ws.Range("R11").Value = 6703432

# Row 12
$ws.Range("A12").Value = 111896652
$ws.Range("B12").Value = 89317
$ws.Range("E12").Value = 3215
$ws.Range("F12").Value = "Rödgul trumpetsvamp"
$ws.Range("G12").Value = "Craterellus lutescens"
$ws.Range("H12").Value = "(Fr.) Fr."
$ws.Range("Q12").Value = 575067
$ws.Range("R12").Value = 6703456

# Row 13
$ws.Range("A13").Value = 111896655
$ws.Range("B13").Value = 89317
$ws.Range("E13").Value = 3215
$ws.Range("F13").Value = "Rödgul trumpetsvamp"
$ws.Range("G13").Value = "Craterellus lutescens"
$ws.Range("H13").Value = "(Fr.) Fr."
$ws.Range("Q13").Value = 575105
$ws.Range("R13").Value = 6703429

# Row 14
$ws.Range("A14").Value = 111896642
$ws.Range("B14").Value = 90466
$ws.Range("Q14").Value = 575014
$ws.Range("R14").Value = 6703387

# Row 15
$ws.Range("A15").Value = 111896634
$ws.Range("B15").Value = 90466
$ws.Range("E15").Value = 4769
$ws.Range("F15").Value = "Svavelriska"
$ws.Range("G15").Value = "Lactarius scrobiculatus"
$ws.Range("H15").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q15").Value = 575048
$ws.Range("R15").Value = 6703452

# Row 16
$ws.Range("A16").Value = 111896636
$ws.Range("B16").Value = 90466
$ws.Range("E16").Value = 4769
$ws.Range("F16").Value = "Svavelriska"
$ws.Range("G16").Value = "Lactarius scrobiculatus"
$ws.Range("H16").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q16").Value = 575109
$ws.Range("R16").Value = 6703418

# Row 17
$ws.Range("A17").Value = 111896635
$ws.Range("B17").Value = 90466
$ws.Range("E17").Value = 4769
$ws.Range("F17").Value = "Svavelriska"
$ws.Range("G17").Value = "Lactarius scrobiculatus"
$ws.Range("H17").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q17").Value = 575037
$ws.Range("R17").Value = 6703389

# Row 18
$ws.Range("A18").Value = 111883983
$ws.Range("B18").Value = 90466
$ws.Range("P18").Value = "Kalkberget (Kalkberget), Gstr"
$ws.Range("Q18").Value = 575058
$ws.Range("R18").Value = 6703446
$ws.Range("AW18").Value = "Patric Engfeldt"
$ws.Range("AX18").Value = "Patric Engfeldt"

# Row 19
$ws.Range("A19").Value = 111896641
$ws.Range("B19").Value = 90466
$ws.Range("E19").Value = 4769
$ws.Range("F19").Value = "Svavelriska"
$ws.Range("G19").Value = "Lactarius scrobiculatus"
$ws.Range("H19").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q19").Value = 575021
$ws.Range("R19").Value = 6703371

# Row 20
$ws.Range("A20").Value = 111896633
$ws.Range("B20").Value = 90466
$ws.Range("Q20").Value = 575100
$ws.Range("R20").Value = 6703444

# Row 21
$ws.Range("A21").Value = 111884471
$ws.Range("B21").Value = 89033
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 3286
$ws.Range("F21").Value = "Flattoppad klubbsvamp"
$ws.Range("G21").Value = "Clavariadelphus truncatus"
$ws.Range("H21").Value = "(Quél.) Donk"
$ws.Range("P21").Value = "Kalkberget (Kalkberget), Gstr"
$ws.Range("Q21").Value = 575021
$ws.Range("R21").Value = 6703397

# Row 22
$ws.Range("A22").Value = 111896653
$ws.Range("B22").Value = 89317
$ws.Range("E22").Value = 3215
$ws.Range("F22").Value = "Rödgul trumpetsvamp"
$ws.Range("G22").Value = "Craterellus lutescens"
$ws.Range("H22").Value = "(Fr.) Fr."
$ws.Range("Q22").Value = 575075
$ws.Range("R22").Value = 6703404

# Row 23
$ws.Range("B23").Value = 90821

# Row 24
$ws.Range("B24").Value = 89100
